# Add customer 345678 to mock data (new row 4 in Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A4: cif_id_mask "345678" -- force text so it isn't coerced to a number,
# then drop back to the unstyled "Normal" style so no extra formatting
# is left on the cell (matches the other data rows, which carry no style).
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "345678"
$ws.Range("A4").Style = "Normal"

# B4: cc_account_open_date "2020-03-15" -- same text-forcing trick so it
# stays a literal date-like string instead of becoming a date serial.
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "2020-03-15"
$ws.Range("B4").Style = "Normal"

# C4: embossed_bin_desc -- plain text, no ambiguity.
$ws.Range("C4").Value = "Visa Signature"

# D4: cc_credit_limit -- numeric value.
$ws.Range("D4").Value = 75000

# E4: cc_account_closed_date -- present but empty (like E2/E3), kept as
# an empty text cell via the quote-prefix trick, then style reset again.
$ws.Range("E4").Value = "'"
$ws.Range("E4").Style = "Normal"
